# Update attendance list for second group:
# Insert a new weekly attendance column (S) right before the existing
# "17.05.2018" column, dated the same as the prior week (10.05.2018),
# and fill in that week's attendance for each student. Also record the
# previously-missing attendance value for row 6 in column R.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column at S - this shifts old S,T,U -> T,U,V (only row 1
# had data there), and the new column S inherits formatting (style) from
# its neighbor, matching column group R:U (width 14).
$ws.Columns("S:S").Insert()

# New date for the inserted week column (same week value as column R,
# matching the committed data).
$ws.Range("S1").Value = $ws.Range("R1").Value()

# Fill in attendance for the new week (column S) per student row.
$ws.Range("S3").Value = $ws.Range("R3").Value()
$ws.Range("S4").Value = $ws.Range("R4").Value()
$ws.Range("S5").Value = $ws.Range("R5").Value()

# Row 6 was missing an entry in column R; add it along with column S.
$ws.Range("R6").Value = $ws.Range("R5").Value()
$ws.Range("S6").Value = $ws.Range("R5").Value()

$ws.Range("S7").Value = $ws.Range("R7").Value()
$ws.Range("S8").Value = $ws.Range("R8").Value()
$ws.Range("S9").Value = $ws.Range("R9").Value()

# Restore/update the active selection to where the user finished editing.
$ws.Range("S9").Select()
